$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Contract Number" value for the new row (A3), which adds a
# new shared string "C037162".
$ws.Range("A3").Value = "C037162"

# Move the active selection to B1 (matches the saved selection state).
$ws.Range("B1").Select()
